# Generate Report for Handback
# Update the Correspond Handoff/Handback Datetime values on the
# per-locale sheets to reflect the new report generation timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-23 13:18:14"
$wsZhCn.Range("H2").Value = "2016-03-23 13:18:37"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-23 13:18:18"
$wsDeDe.Range("H2").Value = "2016-03-23 13:18:44"
